$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1933.3334
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1900
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1900
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2250
# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 9750
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 9750
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9750
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -10718
# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 206.875
$ws.Range("I53").Value = 165.8
$ws.Range("K53").Value = 165.8
$ws.Range("M53").Value = 471.2
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("M116").Value = -558
# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 8124.75
$ws.Range("I127").Value = 7833
$ws.Range("J127").Value = 9000
$ws.Range("K127").Value = 23499
$ws.Range("L127").Value = 27000
$ws.Range("M127").Value = -18539
$ws.Range("N127").Value = -36920
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 621.6667
$ws.Range("I138").Value = 621.6667
$ws.Range("K138").Value = 1865.0001
$ws.Range("M138").Value = 3274.9999

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 1594.6
$ws.Range("I102").Value = 1493.25
$ws.Range("K102").Value = 1493.25
$ws.Range("M102").Value = 128.75
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 9874.875
$ws.Range("I132").Value = 3749.75
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 11249.25
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -8719.25
$ws.Range("N132").Value = -53060

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 3078.75
$ws.Range("I20").Value = 3078.75
$ws.Range("K20").Value = 3078.75
$ws.Range("M20").Value = -2831.75
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 465.2143
$ws.Range("I22").Value = 455.81818
$ws.Range("K22").Value = 455.81818
$ws.Range("M22").Value = -282.81818

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 48 (Leve Item ID 3870)
$ws.Range("H48").Value = 7995
$ws.Range("J48").Value = 7995
$ws.Range("L48").Value = 7995
$ws.Range("N48").Value = -8947
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 4625
$ws.Range("I99").Value = 4250
$ws.Range("K99").Value = 4250
$ws.Range("M99").Value = -2752
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2900
$ws.Range("I105").Value = 2900
$ws.Range("K105").Value = 2900
$ws.Range("M105").Value = -1153
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 991.7778
$ws.Range("I122").Value = 991.7778
$ws.Range("K122").Value = 2975.3334
$ws.Range("M122").Value = -525.3334
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 4625
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 10799.8
$ws.Range("I132").Value = 11999.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 35998.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -33468.5
$ws.Range("N132").Value = -35060
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 5115.0713
$ws.Range("I134").Value = 2688.5
$ws.Range("J134").Value = 8350.5
$ws.Range("K134").Value = 8065.5
$ws.Range("L134").Value = 25051.5
$ws.Range("M134").Value = -5530.5
$ws.Range("N134").Value = -30121.5

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 137
$ws.Range("I46").Value = 132.66667
$ws.Range("J46").Value = 150
$ws.Range("K46").Value = 398.00001
$ws.Range("L46").Value = 450
$ws.Range("M46").Value = -307.00001
$ws.Range("N46").Value = -632
# Row 50 (Leve Item ID 4725)
$ws.Range("H50").Value = 268
$ws.Range("I50").Value = 268
$ws.Range("K50").Value = 804
$ws.Range("M50").Value = -323
# Row 53 (Leve Item ID 4725)
$ws.Range("H53").Value = 268
$ws.Range("I53").Value = 268
$ws.Range("K53").Value = 804
$ws.Range("M53").Value = -323
# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 1492.5714
$ws.Range("I86").Value = 433
$ws.Range("J86").Value = 2287.25
$ws.Range("K86").Value = 1299
$ws.Range("L86").Value = 6861.75
$ws.Range("M86").Value = -113
$ws.Range("N86").Value = -9233.75
# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 1492.5714
$ws.Range("I89").Value = 433
$ws.Range("J89").Value = 2287.25
$ws.Range("K89").Value = 3897
$ws.Range("L89").Value = 20585.25
$ws.Range("M89").Value = 2031
$ws.Range("N89").Value = -32441.25

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2866.3333
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2799.5
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2799.5
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4795.5
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2866.3333
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2799.5
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 13997.5
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -23981.5
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 589.875
$ws.Range("I113").Value = 568.7143
$ws.Range("K113").Value = 568.7143
$ws.Range("M113").Value = 1601.2857

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 13 (Leve Item ID 3546)
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# Row 50 (Leve Item ID 3426)
$ws.Range("H50").Value = 60000
$ws.Range("I50").Value = 60000
$ws.Range("K50").Value = 60000
$ws.Range("M50").Value = -59363
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 80 (Leve Item ID 12027)
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246
# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 2249.5
$ws.Range("I82").Value = 1582.8334
$ws.Range("J82").Value = 3249.5
$ws.Range("K82").Value = 1582.8334
$ws.Range("L82").Value = 3249.5
$ws.Range("M82").Value = -1221.8334
$ws.Range("N82").Value = -3971.5
# Row 83 (Leve Item ID 12027)
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232
# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 2249.5
$ws.Range("I85").Value = 1582.8334
$ws.Range("J85").Value = 3249.5
$ws.Range("K85").Value = 1582.8334
$ws.Range("L85").Value = 3249.5
$ws.Range("M85").Value = -334.8334
$ws.Range("N85").Value = -5745.5
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2936.2727
$ws.Range("I93").Value = 2799.6667
$ws.Range("K93").Value = 2799.6667
$ws.Range("M93").Value = -1551.6667
# Row 99 (Leve Item ID 19636)
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2306
$ws.Range("I100").Value = 1846.75
$ws.Range("J100").Value = 3224.5
$ws.Range("K100").Value = 1846.75
$ws.Range("L100").Value = 3224.5
$ws.Range("M100").Value = -1305.75
$ws.Range("N100").Value = -4306.5
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
# Row 118 (Leve Item ID 26146)
$ws.Range("H118").Value = 24400
$ws.Range("J118").Value = 24400
$ws.Range("L118").Value = 24400
$ws.Range("N118").Value = -27714

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 8 (Leve Item ID 2999)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
# Row 116 (Leve Item ID 26145)
$ws.Range("H116").Value = 40700
$ws.Range("J116").Value = 40700
$ws.Range("L116").Value = 40700
$ws.Range("N116").Value = -49878

